# Daily attendance processing - 2025-10-08 19:17:06
# Swap the order of the two "Recorded By" names in column G so that
# dnasr281@gmail.com is listed first, for the specific rows affected
# by this run (rows where the value is exactly "System, dnasr281@gmail.com"
# or "admin@admin.com, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, admin@admin.com"
    }
}
